$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0.001
$ws.Range("K4").Value = 475
$ws.Range("L4").Value = 0.00095
